$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (Atl. Tembetary - Sportivo Trinidense) and row 32 (Recoleta - 2 de Mayo)
# are removed; the following rows shift up, so former row 33
# (Hamilton - Raith) becomes the new row 31, and the sheet's used range
# shrinks from A1:AS33 to A1:AS31.
$ws.Range("A31:A32").EntireRow.Delete()
